# Add the latest announcement ("*ST海伦" / 300201) to both sheets.
# In each sheet, the new row is inserted immediately above the existing
# last row ("*ST计通" / 300330), pushing that last row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "公告明细" (announcement details)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# The existing last row (26, "*ST计通") needs to move to row 27, so insert
# a new blank row at its current position.
$ws1.Rows.Item(26).Insert()

$ws1.Range("A26").Value = "*ST海伦"

# Force the stock code to be stored as text (like the other codes in the
# column), not auto-converted to a number, then drop back to the default
# "Normal" style so no stray number-format style is left on the cell.
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("B26").Value = "300201"
$ws1.Range("B26").Style = "Normal"

$ws1.Range("C26").Value = "*ST海伦:关于收到《行政处罚及市场禁入事先告知书》的公告"
$ws1.Range("D26").Value = "2023-03-10 00:00:00"
$ws1.Range("E26").Value = "2023-03-09 21:57:56:000"
$ws1.Range("F26").Value = "https://data.eastmoney.com/notices/detail/300201/AN202303091584160746.html"

# ---------------------------------------------------------------------
# Sheet 2: "公告汇总" (announcement summary)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Same pattern: push the existing last row (16, "*ST计通") down to row 17.
$ws2.Rows.Item(16).Insert()

$ws2.Range("A16").Value = "*ST海伦"

$ws2.Range("B16").NumberFormat = "@"
$ws2.Range("B16").Value = "300201"
$ws2.Range("B16").Style = "Normal"

$ws2.Range("C16").Value = 1
